$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 26).Value = 30.04068289318306
$ws.Cells.Item(3, 26).Value = 22.10032906584962
$ws.Cells.Item(4, 26).Value = 15.71643411866143
$ws.Cells.Item(5, 26).Value = 10.6161421138269
$ws.Cells.Item(6, 26).Value = 6.62608620842262
$ws.Cells.Item(7, 26).Value = 3.632006118520721
$ws.Cells.Item(8, 26).Value = 1.58401478028587
$ws.Cells.Item(9, 26).Value = 0.3819458071278916
$ws.Cells.Item(10, 26).Value = 0.000375000585933094
$ws.Cells.Item(11, 26).Value = 0.6686231538410539
$ws.Cells.Item(12, 26).Value = 2.633054056944971
$ws.Cells.Item(13, 5).Value = 180
$ws.Cells.Item(13, 6).Value = 5.217597167510036
$ws.Cells.Item(13, 8).Value = 175
$ws.Cells.Item(13, 9).Value = 3.810499613777495
$ws.Cells.Item(13, 11).Value = 43.56059998523794
$ws.Cells.Item(13, 12).Value = 10
$ws.Cells.Item(13, 13).Value = 14.47751224073088
$ws.Cells.Item(13, 14).Value = 10.327955592437
$ws.Cells.Item(13, 17).Value = 36.18950038622251
$ws.Cells.Item(13, 18).Value = 102570.1405098137
$ws.Cells.Item(13, 19).Value = 2500
$ws.Cells.Item(13, 20).Value = 2500
$ws.Cells.Item(13, 21).Value = 25000
$ws.Cells.Item(13, 22).Value = 175
$ws.Cells.Item(13, 26).Value = 5.217597167510036
$ws.Cells.Item(13, 32).Value = 68.5
$ws.Cells.Item(13, 33).Value = 64.68950038622251
$ws.Cells.Item(13, 34).Value = 4036.624824100284
$ws.Cells.Item(14, 2).Value = 180
$ws.Cells.Item(14, 3).Value = 5.217597167510036
$ws.Cells.Item(14, 5).Value = 190
$ws.Cells.Item(14, 6).Value = 8.757022693565062
$ws.Cells.Item(14, 8).Value = 185
$ws.Cells.Item(14, 9).Value = 6.862915010152392
$ws.Cells.Item(14, 11).Value = 45.12351630185047
$ws.Cells.Item(14, 12).Value = 10
$ws.Cells.Item(14, 13).Value = 19.47122071361051
$ws.Cells.Item(14, 14).Value = 10.60660172297659
$ws.Cells.Item(14, 17).Value = 33.13708498984761
$ws.Cells.Item(14, 18).Value = 98540.95218659885
$ws.Cells.Item(14, 19).Value = 2500
$ws.Cells.Item(14, 20).Value = 2500
$ws.Cells.Item(14, 21).Value = 25000
$ws.Cells.Item(14, 22).Value = 185
$ws.Cells.Item(14, 26).Value = 8.757022693565062
$ws.Cells.Item(14, 32).Value = 71.07142857142857
$ws.Cells.Item(14, 33).Value = 64.20851356127618
$ws.Cells.Item(14, 34).Value = 4006.611246223633
$ws.Cells.Item(15, 2).Value = 190
$ws.Cells.Item(15, 3).Value = 8.757022693565062
$ws.Cells.Item(15, 5).Value = 200
$ws.Cells.Item(15, 6).Value = 13.3463549614923
$ws.Cells.Item(15, 8).Value = 195
$ws.Cells.Item(15, 9).Value = 10.91287885364285
$ws.Cells.Item(15, 11).Value = 47.20210304931079
$ws.Cells.Item(15, 12).Value = 10
$ws.Cells.Item(15, 13).Value = 24.62431846249069
$ws.Cells.Item(15, 14).Value = 11.00038197404725
$ws.Cells.Item(15, 17).Value = 29.08712114635715
$ws.Cells.Item(15, 18).Value = 93194.99991319144
$ws.Cells.Item(15, 19).Value = 2500
$ws.Cells.Item(15, 20).Value = 2500
$ws.Cells.Item(15, 21).Value = 25000
$ws.Cells.Item(15, 22).Value = 195
$ws.Cells.Item(15, 26).Value = 13.3463549614923
$ws.Cells.Item(15, 32).Value = 73.21428571428571
$ws.Cells.Item(15, 33).Value = 62.30140686064286
$ws.Cells.Item(15, 34).Value = 3887.607788104114
$ws.Cells.Item(16, 2).Value = 200
$ws.Cells.Item(16, 3).Value = 13.3463549614923
$ws.Cells.Item(16, 5).Value = 210
$ws.Cells.Item(16, 6).Value = 19.12879499083992
$ws.Cells.Item(16, 8).Value = 205
$ws.Cells.Item(16, 9).Value = 16.07695154586736
$ws.Cells.Item(16, 11).Value = 49.86211460811401
$ws.Cells.Item(16, 12).Value = 10
$ws.Cells.Item(16, 13).Value = 30.00000015314095
$ws.Cells.Item(16, 14).Value = 11.54700540161127
$ws.Cells.Item(16, 17).Value = 23.92304845413264
$ws.Cells.Item(16, 18).Value = 86378.42395945509
$ws.Cells.Item(16, 19).Value = 2500
$ws.Cells.Item(16, 20).Value = 2500
$ws.Cells.Item(16, 21).Value = 25000
$ws.Cells.Item(16, 22).Value = 205
$ws.Cells.Item(16, 26).Value = 19.12879499083992
$ws.Cells.Item(16, 32).Value = 75.35714285714286
$ws.Cells.Item(16, 33).Value = 59.2801913112755
$ws.Cells.Item(16, 34).Value = 3699.083937823591
$ws.Cells.Item(17, 2).Value = 210
$ws.Cells.Item(17, 3).Value = 19.12879499083992
$ws.Cells.Item(17, 5).Value = 220
$ws.Cells.Item(17, 6).Value = 26.32503002402403
$ws.Cells.Item(17, 8).Value = 215
$ws.Cells.Item(17, 9).Value = 22.53205655191036
$ws.Cells.Item(17, 11).Value = 53.20633696566401
$ws.Cells.Item(17, 12).Value = 10
$ws.Cells.Item(17, 13).Value = 35.68533492923775
$ws.Cells.Item(17, 14).Value = 12.31174025844619
$ws.Cells.Item(17, 17).Value = 17.46794344808964
$ws.Cells.Item(17, 18).Value = 77857.68535147833
$ws.Cells.Item(17, 19).Value = 2500
$ws.Cells.Item(17, 20).Value = 2500
$ws.Cells.Item(17, 21).Value = 25000
$ws.Cells.Item(17, 22).Value = 215
$ws.Cells.Item(17, 26).Value = 26.32503002402403
$ws.Cells.Item(17, 32).Value = 76.56756756756756
$ws.Cells.Item(17, 33).Value = 54.03551101565721
$ws.Cells.Item(17, 34).Value = 3371.815887377009
$ws.Cells.Item(18, 2).Value = 220
$ws.Cells.Item(18, 3).Value = 26.32503002402403
$ws.Cells.Item(18, 5).Value = 228.9543811126312
$ws.Cells.Item(18, 6).Value = 34.258167199464
$ws.Cells.Item(18, 8).Value = 224.4771905563156
$ws.Cells.Item(18, 9).Value = 30.09240198250708
$ws.Cells.Item(18, 11).Value = 57.16039905094622
$ws.Cells.Item(18, 12).Value = 8.954381112631211
$ws.Cells.Item(18, 13).Value = 41.47627967182778
$ws.Cells.Item(18, 14).Value = 11.95144528831471
$ws.Cells.Item(18, 17).Value = 9.907598017492916
$ws.Cells.Item(18, 18).Value = 60780.57442705733
$ws.Cells.Item(18, 22).Value = 224.4771905563156
$ws.Cells.Item(18, 26).Value = 34.258167199464
$ws.Cells.Item(18, 32).Value = 77.33598842348505
$ws.Cells.Item(18, 33).Value = 47.24358644097796
$ws.Cells.Item(18, 34).Value = 2947.999793917025
$ws.Cells.Item(19, 2).Value = 228.9543811126312
$ws.Cells.Item(19, 3).Value = 34.258167199464
$ws.Cells.Item(19, 5).Value = 237.9087622252624
$ws.Cells.Item(19, 6).Value = 44.05290063623465
$ws.Cells.Item(19, 8).Value = 233.4315716689468
$ws.Cells.Item(19, 9).Value = 38.88368146815391
$ws.Cells.Item(19, 11).Value = 61.8357956447819
$ws.Cells.Item(19, 12).Value = 8.954381112631239
$ws.Cells.Item(19, 13).Value = 47.47053161366336
$ws.Cells.Item(19, 14).Value = 13.24672720197736
$ws.Cells.Item(19, 16).Value = 24
$ws.Cells.Item(19, 17).Value = 1.11631853184609
$ws.Cells.Item(19, 18).Value = 50389.47278541525
$ws.Cells.Item(19, 22).Value = 233.4315716689468
$ws.Cells.Item(19, 26).Value = 44.05290063623465
$ws.Cells.Item(19, 32).Value = 78.06201932450921
$ws.Cells.Item(19, 33).Value = 39.1783378563553
$ws.Cells.Item(19, 34).Value = 2444.728282236571
$ws.Cells.Item(20, 2).Value = 237.9087622252624
$ws.Cells.Item(20, 3).Value = 44.05290063623465
$ws.Cells.Item(20, 5).Value = 246.8631433378937
$ws.Cells.Item(20, 6).Value = 56.56578187347365
$ws.Cells.Item(20, 8).Value = 242.3859527815781
$ws.Cells.Item(20, 9).Value = 49.88597714562194
$ws.Cells.Item(20, 11).Value = 67.27170458889177
$ws.Cells.Item(20, 12).Value = 8.954381112631239
$ws.Cells.Item(20, 13).Value = 54.24761044425681
$ws.Cells.Item(20, 14).Value = 15.32540440350763
$ws.Cells.Item(20, 16).Value = 14.11402285437806
$ws.Cells.Item(20, 18).Value = 38447.27165329985
$ws.Cells.Item(20, 22).Value = 242.3859527815781
$ws.Cells.Item(20, 26).Value = 56.56578187347365
$ws.Cells.Item(20, 32).Value = 78.78805022553335
$ws.Cells.Item(20, 33).Value = 28.90207307991142
$ws.Cells.Item(20, 34).Value = 1803.489360186472
$ws.Cells.Item(21, 2).Value = 246.8631433378937
$ws.Cells.Item(21, 3).Value = 56.56578187347365
$ws.Cells.Item(21, 5).Value = 255.8175244505249
$ws.Cells.Item(21, 6).Value = 73.96222991219801
$ws.Cells.Item(21, 8).Value = 251.3403338942093
$ws.Cells.Item(21, 9).Value = 64.39664230221271
$ws.Cells.Item(21, 11).Value = 74.19832115110636
$ws.Cells.Item(21, 12).Value = 8.954381112631211
$ws.Cells.Item(21, 13).Value = 62.3957812014332
$ws.Cells.Item(21, 14).Value = 19.32483651285239
$ws.Cells.Item(21, 15).Value = 19.60335769778729
$ws.Cells.Item(21, 18).Value = 22819.67166871862
$ws.Cells.Item(21, 22).Value = 251.3403338942093
$ws.Cells.Item(21, 25).Value = 3120
$ws.Cells.Item(21, 26).Value = 77.29556324553134
$ws.Cells.Item(21, 32).Value = 79.08453778525612
$ws.Cells.Item(21, 33).Value = 14.68789548304341
$ws.Cells.Item(21, 34).Value = 916.5246781419085
